$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 86, shifting existing rows 86:201 down to 87:202.
# This preserves the formatting (including the date style on column D) of the
# row that was previously at 86, since Excel's row insert copies formats from
# the row above by default behavior mirrors the row being pushed down.
$ws.Rows.Item(86).Insert()

# Populate the newly inserted row 86 with the new record's data.
$ws.Cells.Item(86, 1).Value = 9
$ws.Cells.Item(86, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(86, 3).Value = "Metropolitana"
$ws.Cells.Item(86, 4).Value = 44665
$ws.Cells.Item(86, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(86, 5).Value = 13
$ws.Cells.Item(86, 6).Value = "Fruta"
$ws.Cells.Item(86, 7).Value = 100101
$ws.Cells.Item(86, 8).Value = "Berries"
$ws.Cells.Item(86, 9).Value = 100101001
$ws.Cells.Item(86, 10).Value = "Arándano (blue)"
$ws.Cells.Item(86, 11).Value = "Sin especificar"
$ws.Cells.Item(86, 12).Value = "Primera"
$ws.Cells.Item(86, 13).Value = 530
$ws.Cells.Item(86, 14).Value = 4500
$ws.Cells.Item(86, 15).Value = 5000
$ws.Cells.Item(86, 16).Value = 4736
$ws.Cells.Item(86, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(86, 18).Value = "Provincia de Linares"
$ws.Cells.Item(86, 19).Value = 2368
$ws.Cells.Item(86, 20).Value = 2
